$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the department name (missing trailing "a")
$ws.Range("B26").Value = "Archipiélago de San Andrés, Providencia y Santa Catalina"

# Move the active selection to B26 (last edited cell)
$ws.Range("B26").Select()
